$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Step 1: Insert two new columns before column D (old D shifts to new F, etc.)
$ws.Range("D1:E1").EntireColumn.Insert()

# Step 2: Copy number formats/styles from the (shifted) old-D column (now F) into new D,E columns
$ws.Range("F7:F102").Copy()
$ws.Range("D7:D102").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E7:E102").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Step 3: Populate new D/E column values (two new quarterly periods) for each data row
$ws.Cells.Item(7, 4).Value = 43465
$ws.Cells.Item(7, 5).Value = 43373
$ws.Cells.Item(8, 4).Value = 3195000
$ws.Cells.Item(8, 5).Value = 2931000
$ws.Cells.Item(9, 4).Value = 1665000
$ws.Cells.Item(9, 5).Value = 1945000
$ws.Cells.Item(10, 4).Value = 1530000
$ws.Cells.Item(10, 5).Value = 986000
$ws.Cells.Item(12, 4).Value = "NA"
$ws.Cells.Item(12, 5).Value = "NA"
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(14, 4).Value = -706000
$ws.Cells.Item(14, 5).Value = 4000
$ws.Cells.Item(15, 4).Value = 391000
$ws.Cells.Item(15, 5).Value = 380000
$ws.Cells.Item(17, 4).Value = 1470000
$ws.Cells.Item(17, 5).Value = 2460000
$ws.Cells.Item(18, 4).Value = 1725000
$ws.Cells.Item(18, 5).Value = 471000
$ws.Cells.Item(20, 4).Value = -424000
$ws.Cells.Item(20, 5).Value = 188000
$ws.Cells.Item(21, 4).Value = 1692000
$ws.Cells.Item(21, 5).Value = 1039000
$ws.Cells.Item(22, 4).Value = 240000
$ws.Cells.Item(22, 5).Value = 232000
$ws.Cells.Item(23, 4).Value = 1061000
$ws.Cells.Item(23, 5).Value = 427000
$ws.Cells.Item(24, 4).Value = 223000
$ws.Cells.Item(24, 5).Value = 167000
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(26, 4).Value = 838000
$ws.Cells.Item(26, 5).Value = 260000
$ws.Cells.Item(27, 4).Value = 864000
$ws.Cells.Item(27, 5).Value = 274000
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(29, 4).Value = "NA"
$ws.Cells.Item(29, 5).Value = "NA"
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(32, 4).Value = 424000
$ws.Cells.Item(32, 5).Value = -188000
$ws.Cells.Item(33, 4).Value = 864000
$ws.Cells.Item(33, 5).Value = 274000
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(35, 4).Value = 864000
$ws.Cells.Item(35, 5).Value = 274000
$ws.Cells.Item(38, 4).Value = 43465
$ws.Cells.Item(38, 5).Value = 43373
$ws.Cells.Item(41, 4).Value = 190000
$ws.Cells.Item(41, 5).Value = 212000
$ws.Cells.Item(42, 4).Value = 0
$ws.Cells.Item(42, 5).Value = 0
$ws.Cells.Item(43, 4).Value = 1918000
$ws.Cells.Item(43, 5).Value = 1762000
$ws.Cells.Item(44, 4).Value = 296000
$ws.Cells.Item(44, 5).Value = 345000
$ws.Cells.Item(45, 4).Value = 1241000
$ws.Cells.Item(45, 5).Value = 2726000
$ws.Cells.Item(46, 4).Value = 3645000
$ws.Cells.Item(46, 5).Value = 5045000
$ws.Cells.Item(47, 4).Value = 14112000
$ws.Cells.Item(47, 5).Value = 14312000
$ws.Cells.Item(48, 4).Value = 36796000
$ws.Cells.Item(48, 5).Value = 35498000
$ws.Cells.Item(49, 4).Value = 2645000
$ws.Cells.Item(49, 5).Value = 2592000
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(50, 5).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(51, 5).Value = 0
$ws.Cells.Item(52, 4).Value = 3440000
$ws.Cells.Item(52, 5).Value = 3158000
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(53, 5).Value = 0
$ws.Cells.Item(54, 4).Value = 60638000
$ws.Cells.Item(54, 5).Value = 60605000
$ws.Cells.Item(57, 4).Value = 1474000
$ws.Cells.Item(57, 5).Value = 1375000
$ws.Cells.Item(58, 4).Value = 3752000
$ws.Cells.Item(58, 5).Value = 4361000
$ws.Cells.Item(59, 4).Value = 2297000
$ws.Cells.Item(59, 5).Value = 2755000
$ws.Cells.Item(60, 4).Value = 7523000
$ws.Cells.Item(60, 5).Value = 8491000
$ws.Cells.Item(61, 4).Value = 21611000
$ws.Cells.Item(61, 5).Value = 21335000
$ws.Cells.Item(62, 4).Value = 12256000
$ws.Cells.Item(62, 5).Value = 11668000
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(63, 5).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(64, 5).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(65, 5).Value = 0
$ws.Cells.Item(66, 4).Value = 43480000
$ws.Cells.Item(66, 5).Value = 43968000
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(68, 5).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(69, 5).Value = 0
$ws.Cells.Item(70, 4).Value = 2278000
$ws.Cells.Item(70, 5).Value = 2279000
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(71, 5).Value = 0
$ws.Cells.Item(72, 4).Value = 10104000
$ws.Cells.Item(72, 5).Value = 9485000
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(73, 5).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(74, 5).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(75, 5).Value = 0
$ws.Cells.Item(76, 4).Value = 14880000
$ws.Cells.Item(76, 5).Value = 14358000
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(77, 5).Value = 0
$ws.Cells.Item(80, 4).Value = 43465
$ws.Cells.Item(80, 5).Value = 43373
$ws.Cells.Item(81, 4).Value = 864000
$ws.Cells.Item(81, 5).Value = 274000
$ws.Cells.Item(83, 4).Value = 391000
$ws.Cells.Item(83, 5).Value = 380000
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(84, 5).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(85, 5).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(86, 5).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(87, 5).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(88, 5).Value = 0
$ws.Cells.Item(89, 4).Value = 856000
$ws.Cells.Item(89, 5).Value = 962000
$ws.Cells.Item(91, 4).Value = -969000
$ws.Cells.Item(91, 5).Value = -874000
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(92, 5).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(93, 5).Value = 0
$ws.Cells.Item(94, 4).Value = 147000
$ws.Cells.Item(94, 5).Value = -964000
$ws.Cells.Item(96, 4).Value = -232000
$ws.Cells.Item(96, 5).Value = -229000
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(97, 5).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(98, 5).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(99, 5).Value = 0
$ws.Cells.Item(100, 4).Value = -1039000
$ws.Cells.Item(100, 5).Value = -32000
$ws.Cells.Item(101, 4).Value = -6000
$ws.Cells.Item(101, 5).Value = -5000
$ws.Cells.Item(102, 4).Value = -42000
$ws.Cells.Item(102, 5).Value = -39000
$ws.Range("D:E").EntireColumn.AutoFit()
